# SwaadSutra Daily Orders - 2026-01-13 - new order (#9) received
# Order #9 (Sagar Borse, Wheat Chapati x1) is added to the top of the Daily
# Orders log, pushing every previous order down by one row. Summary and
# Items Breakdown roll-ups are updated to match.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Daily Orders ---
$ws1 = $wb.Worksheets.Item("Daily Orders")

# Duplicate the current top row and insert the copy above it. This keeps the
# new row's cell typing (plain numbers vs. text-stored phone/date strings)
# consistent with the rest of the sheet instead of letting Excel's
# auto-detection turn things like the phone number or collection date into
# numeric/date values.
$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 9
$ws1.Cells.Item(2, 2).Value = "2026-01-13 19:05"
$ws1.Cells.Item(2, 3).Value = "Sagar Borse"
$ws1.Cells.Item(2, 4).Value = "A1608"
$ws1.Cells.Item(2, 5).Value = "'7588930329"
$ws1.Cells.Item(2, 6).Value = "Wheat Chapati x1"
$ws1.Cells.Item(2, 7).Value = 15
$ws1.Cells.Item(2, 8).Value = "NEW"
$ws1.Cells.Item(2, 9).Value = "PENDING"
$ws1.Cells.Item(2, 10).Value = "'2026-01-26"
$ws1.Cells.Item(2, 11).Value = "10:35"

# --- Sheet 2: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 1).Value = 9    # Total Orders: 8 -> 9
$ws2.Cells.Item(2, 2).Value = 8    # New: 7 -> 8
$ws2.Cells.Item(2, 7).Value = 230  # Total Revenue: 215 -> 230

# --- Sheet 3: Items Breakdown ---
$ws3 = $wb.Worksheets.Item("Items Breakdown")
$ws3.Cells.Item(2, 2).Value = 4    # Wheat Chapati quantity: 3 -> 4
$ws3.Cells.Item(2, 3).Value = 60   # Wheat Chapati revenue: 45 -> 60
